$d = $word.ActiveDocument
$t = $d.Tables(1)

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Set-CellXml($row, $col, $pAttrs, $innerXml) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Delete()
    $cell2 = $t.Cell($row, $col)
    $xml = '<w:p ' + $wns + ' ' + $pAttrs + '>' + $innerXml + '</w:p>'
    $cell2.Range.InsertXML($xml)
}

# --- Merge split runs / remove proofErr spans (de-duplicated text) ---

# E2000 row, "Mensaje" column: "La solicitud con el id {requestID} no existe."
Set-CellXml 3 4 'w14:paraId="4A3697B9" w14:textId="708EEB7B" w:rsidR="00DF3010" w:rsidRPr="00F85C11" w:rsidRDefault="00F85C11"' `
    ('<w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr>' + `
     '<w:r w:rsidRPr="00F85C11"><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>La solicitud con el id {requestID} no existe.</w:t></w:r>')

# E3000 row, "Mensaje" column: 'El usuario con el id \"{userId}\" no existe.'
Set-CellXml 4 4 'w14:paraId="05336E08" w14:textId="2ECB02F9" w:rsidR="0013115B" w:rsidRPr="00122D48" w:rsidRDefault="00122D48"' `
    ('<w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr>' + `
     '<w:r w:rsidRPr="00122D48"><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>El usuario con el id \"{userId}\" no existe.</w:t></w:r>')

# E3000 row, "Escenario" column: "Al intentar retornar un usuario mediante un id."
Set-CellXml 4 5 'w14:paraId="2110E4A0" w14:textId="162E7A79" w:rsidR="0013115B" w:rsidRDefault="00CB4701"' `
    ('<w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr>' + `
     '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Al intentar retornar un usuario mediante un id.</w:t></w:r>')

# E3002 row, "Titulo" column: "Formato de token incorrecto"
Set-CellXml 6 2 'w14:paraId="7C1CF9E8" w14:textId="2AB2B25D" w:rsidR="00C55DE0" w:rsidRDefault="00F308D8"' `
    '<w:r><w:t>Formato de token incorrecto</w:t></w:r>'

# E3003 row, "Titulo" column: "Token invalido"
Set-CellXml 7 2 'w14:paraId="0A0A5003" w14:textId="612A01B6" w:rsidR="00564895" w:rsidRDefault="00564895"' `
    '<w:r><w:t>Token invalido</w:t></w:r>'

# E3003 row, "Escenario" column: "Al intentar verificar un correo electrónico con un token invalido."
Set-CellXml 7 5 'w14:paraId="7760833A" w14:textId="4192C4B3" w:rsidR="00564895" w:rsidRDefault="00564895"' `
    ('<w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr>' + `
     '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Al intentar verificar un correo electrónico con un token invalido.</w:t></w:r>')

# --- Add new row: E3004 / Credenciales invalidas / 401 ---

$newRow = $t.Rows.Add()
$r = $t.Rows.Count

$t.Cell($r, 1).Range.Text = "E3004"
$t.Cell($r, 2).Range.Text = "Credenciales invalidas"
$t.Cell($r, 3).Range.Text = "Incorrect username or password."
$t.Cell($r, 4).Range.Text = "Usuario y/o contraseña incorrectos."
$t.Cell($r, 5).Range.Text = "Al intentar iniciar sesión con credenciales incorrectas."
$t.Cell($r, 6).Range.Text = "401"
